# Update the 10-year historical return values for rows 2-5 (open/high/low/close),
# columns F:AE, per the refreshed DataCamp "ten_yr" dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"6.0999999999999999E-2"
$ws.Range("G2").Value = [double]"2.7E-2"
$ws.Range("H2").Value = [double]"-5.8999999999999997E-2"
$ws.Range("I2").Value = [double]"-4.5999999999999999E-2"
$ws.Range("J2").Value = [double]"1.4E-2"
$ws.Range("K2").Value = [double]"-2.7E-2"
$ws.Range("L2").Value = [double]"-0.115"
$ws.Range("M2").Value = [double]"2.3E-2"
$ws.Range("N2").Value = [double]"-9.4E-2"
$ws.Range("O2").Value = [double]"-2.5000000000000001E-2"
$ws.Range("P2").Value = [double]"-2.1000000000000001E-2"
$ws.Range("Q2").Value = [double]"7.5999999999999998E-2"
$ws.Range("R2").Value = [double]"7.2999999999999995E-2"
$ws.Range("S2").Value = [double]"-2.1000000000000001E-2"
$ws.Range("T2").Value = [double]"0"
$ws.Range("U2").Value = [double]"-2.4E-2"
$ws.Range("V2").Value = [double]"-2.5999999999999999E-2"
$ws.Range("W2").Value = [double]"0.05"
$ws.Range("X2").Value = [double]"-0.26800000000000002"
$ws.Range("Y2").Value = [double]"-0.23599999999999999"
$ws.Range("Z2").Value = [double]"0.28599999999999998"
$ws.Range("AA2").Value = [double]"4.5999999999999999E-2"
$ws.Range("AB2").Value = [double]"-0.10299999999999999"
$ws.Range("AC2").Value = [double]"0.191"
$ws.Range("AD2").Value = [double]"0.107"
$ws.Range("AE2").Value = [double]"2.4E-2"

# Row 3
$ws.Range("F3").Value = [double]"0.08"
$ws.Range("G3").Value = [double]"-2.1999999999999999E-2"
$ws.Range("H3").Value = [double]"-0.06"
$ws.Range("I3").Value = [double]"-3.7999999999999999E-2"
$ws.Range("J3").Value = [double]"4.0000000000000001E-3"
$ws.Range("K3").Value = [double]"-5.6000000000000001E-2"
$ws.Range("L3").Value = [double]"-3.9E-2"
$ws.Range("M3").Value = [double]"-5.2999999999999999E-2"
$ws.Range("N3").Value = [double]"-2.3E-2"
$ws.Range("O3").Value = [double]"-6.5000000000000002E-2"
$ws.Range("P3").Value = [double]"0.05"
$ws.Range("Q3").Value = [double]"6.5000000000000002E-2"
$ws.Range("R3").Value = [double]"4.4999999999999998E-2"
$ws.Range("S3").Value = [double]"-3.5000000000000003E-2"
$ws.Range("T3").Value = [double]"-0.02"
$ws.Range("U3").Value = [double]"-4.4999999999999998E-2"
$ws.Range("V3").Value = [double]"5.1999999999999998E-2"
$ws.Range("W3").Value = [double]"-3.5000000000000003E-2"
$ws.Range("X3").Value = [double]"-0.27200000000000002"
$ws.Range("Y3").Value = [double]"-1.0999999999999999E-2"
$ws.Range("Z3").Value = [double]"6.9000000000000006E-2"
$ws.Range("AA3").Value = [double]"-4.0000000000000001E-3"
$ws.Range("AB3").Value = [double]"4.1000000000000002E-2"
$ws.Range("AC3").Value = [double]"0.187"
$ws.Range("AD3").Value = [double]"6.8000000000000005E-2"
$ws.Range("AE3").Value = [double]"-6.2E-2"

# Row 4
$ws.Range("F4").Value = [double]"5.8999999999999997E-2"
$ws.Range("G4").Value = [double]"-2.7E-2"
$ws.Range("H4").Value = [double]"-5.1999999999999998E-2"
$ws.Range("I4").Value = [double]"-4.2999999999999997E-2"
$ws.Range("J4").Value = [double]"3.0000000000000001E-3"
$ws.Range("K4").Value = [double]"-0.108"
$ws.Range("L4").Value = [double]"-1E-3"
$ws.Range("M4").Value = [double]"-0.14599999999999999"
$ws.Range("N4").Value = [double]"7.4999999999999997E-2"
$ws.Range("O4").Value = [double]"-6.7000000000000004E-2"
$ws.Range("P4").Value = [double]"4.2999999999999997E-2"
$ws.Range("Q4").Value = [double]"7.2999999999999995E-2"
$ws.Range("R4").Value = [double]"4.7E-2"
$ws.Range("S4").Value = [double]"-2.3E-2"
$ws.Range("T4").Value = [double]"-2E-3"
$ws.Range("U4").Value = [double]"-0.13600000000000001"
$ws.Range("V4").Value = [double]"4.5999999999999999E-2"
$ws.Range("W4").Value = [double]"-0.13800000000000001"
$ws.Range("X4").Value = [double]"-0.30399999999999999"
$ws.Range("Y4").Value = [double]"5.8999999999999997E-2"
$ws.Range("Z4").Value = [double]"0.217"
$ws.Range("AA4").Value = [double]"-6.2E-2"
$ws.Range("AB4").Value = [double]"6.9000000000000006E-2"
$ws.Range("AC4").Value = [double]"0.16800000000000001"
$ws.Range("AD4").Value = [double]"0.123"
$ws.Range("AE4").Value = [double]"-5.5E-2"

# Row 5
$ws.Range("F5").Value = [double]"2.9000000000000001E-2"
$ws.Range("G5").Value = [double]"-5.1999999999999998E-2"
$ws.Range("H5").Value = [double]"-4.9000000000000002E-2"
$ws.Range("I5").Value = [double]"8.9999999999999993E-3"
$ws.Range("J5").Value = [double]"-2.3E-2"
$ws.Range("K5").Value = [double]"-0.112"
$ws.Range("L5").Value = [double]"1.6E-2"
$ws.Range("M5").Value = [double]"-9.8000000000000004E-2"
$ws.Range("N5").Value = [double]"-2.9000000000000001E-2"
$ws.Range("O5").Value = [double]"-2.9000000000000001E-2"
$ws.Range("P5").Value = [double]"9.5000000000000001E-2"
$ws.Range("Q5").Value = [double]"7.5999999999999998E-2"
$ws.Range("R5").Value = [double]"-1.7000000000000001E-2"
$ws.Range("S5").Value = [double]"0"
$ws.Range("T5").Value = [double]"-4.2000000000000003E-2"
$ws.Range("U5").Value = [double]"4.0000000000000001E-3"
$ws.Range("V5").Value = [double]"3.6999999999999998E-2"
$ws.Range("W5").Value = [double]"-0.255"
$ws.Range("X5").Value = [double]"-0.24099999999999999"
$ws.Range("Y5").Value = [double]"0.26700000000000002"
$ws.Range("Z5").Value = [double]"6.9000000000000006E-2"
$ws.Range("AA5").Value = [double]"-0.11700000000000001"
$ws.Range("AB5").Value = [double]"0.16400000000000001"
$ws.Range("AC5").Value = [double]"0.109"
$ws.Range("AD5").Value = [double]"1.7000000000000001E-2"
$ws.Range("AE5").Value = [double]"-6.0000000000000001E-3"

# Row 1 date headers: switch the custom "yyyy-mm-dd" number format to the
# built-in short-date format.
$ws.Range("B1").NumberFormat = "mm-dd-yy"
$ws.Range("B1").Copy()
$ws.Range("C1:AJ1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the last-used selection.
$ws.Range("D11").Select()
